# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "fdd33410-614b-4b30-8a6c-e22f865d3153" item in both
# the zh-cn and de-de locale sheets (and the rolled-up Overview sheet).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Row 3 on every sheet corresponds to the fdd33410-... file; update its
# status from "Ready for handoff" to "Handback transform failed".
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch error detail for each
# locale in column P ("Error Detail") of row 3.
$zhCnError = "Handback file name: jttajqbv.yqn is different with handoff file name: fdd33410-614b-4b30-8a6c-e22f865d3153.90deb617c8657eff29ba21828b3c30b77be85ad9.zh-cn."
$deDeError = "Handback file name: jttajqbv.yqn is different with handoff file name: fdd33410-614b-4b30-8a6c-e22f865d3153.90deb617c8657eff29ba21828b3c30b77be85ad9.de-de."

$wsZhCn.Range("P3").Value = $zhCnError
$wsDeDe.Range("P3").Value = $deDeError

# Widen the "Error Detail" column (column 16 / P) so the longer error
# messages are readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
